$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose new values look like numbers need to be forced to stay text
# (matching the original inlineStr/text representation), then have their
# temporary "Text" number format cleared so no extra style is left behind.
$textForceCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D26", "D27", "D29", "D30", "D31", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated Price (D) and Volume(1h) (E) values, and the Coin name (B)
# / Link (C) swaps, exactly as described by the diff.
$ws.Range("D2").Value = '28.188.96'
$ws.Range("E2").Value = '  -1.62%  '
$ws.Range("D3").Value = '1.804.64'
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '316.43'
$ws.Range("E5").Value = '  +0.98%  '
$ws.Range("D6").Value = '1.0000'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '0.5377'
$ws.Range("E7").Value = '  +1.49%  '
$ws.Range("D8").Value = '0.3787'
$ws.Range("E8").Value = '  +0.55%  '
$ws.Range("D9").Value = '0.07475'
$ws.Range("E9").Value = '  -0.77%  '
$ws.Range("D10").Value = '42.02'
$ws.Range("E10").Value = '  -1.03%  '
$ws.Range("D11").Value = '1.098'
$ws.Range("E11").Value = '  -2.06%  '
$ws.Range("D12").Value = '0.9993'
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").Value = '6.207'
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("D14").Value = '20.52'
$ws.Range("E14").Value = '  -3.09%  '
$ws.Range("D15").Value = '7.387'
$ws.Range("E15").Value = '  -0.86%  '
$ws.Range("D16").Value = '1.804.48'
$ws.Range("E16").Value = '  +0.71%  '
$ws.Range("D17").Value = '89.88'
$ws.Range("E17").Value = '  -0.63%  '
$ws.Range("D18").Value = '0.00001065'
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").Value = '0.06497'
$ws.Range("E19").Value = '  +0.81%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '17.39'
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").Value = '5.931'
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").Value = '28.206.94'
$ws.Range("E23").Value = '  -1.57%  '
$ws.Range("D24").Value = '11.21'
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("E25").Value = '  -0.29%  '
$ws.Range("D26").Value = '156.14'
$ws.Range("E26").Value = '  -3.01%  '
$ws.Range("D27").Value = '20.53'
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("D28").Value = '2.010.71'
$ws.Range("E28").Value = '  +0.65%  '
$ws.Range("D29").Value = '2.339'
$ws.Range("E29").Value = '  -2.76%  '
$ws.Range("D30").Value = '122.06'
$ws.Range("E30").Value = '  -1.21%  '
$ws.Range("D31").Value = '1.135'
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("E32").Value = '  +9.02%  '
$ws.Range("D33").Value = '3.668'
$ws.Range("E33").Value = '  +0.12%  '
$ws.Range("D34").Value = '5.601'
$ws.Range("E34").Value = '  -2.27%  '
$ws.Range("E35").Value = '  +7.29%  '
$ws.Range("D36").Value = '0.2225'
$ws.Range("E36").Value = '  -3.18%  '
$ws.Range("D37").Value = '0.02303'
$ws.Range("E37").Value = '  -1.02%  '
$ws.Range("D38").Value = '5.093'
$ws.Range("E38").Value = '  +0.30%  '
$ws.Range("D39").Value = '8.469'
$ws.Range("E39").Value = '  -3.81%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '11.15'
$ws.Range("E40").Value = '  -3.07%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.6177'
$ws.Range("E41").Value = '  -2.30%  '
$ws.Range("B42").Value = 'WEMIXTOKEN'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").Value = '1.430'
$ws.Range("E42").Value = '  +2.47%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '1.175'
$ws.Range("E43").Value = '  -2.51%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '13.38'
$ws.Range("E45").Value = '  -1.37%  '
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = '3.682'
$ws.Range("E46").Value = '  +0.46%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '0.5768'
$ws.Range("E47").Value = '  -2.64%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '125.44'
$ws.Range("E48").Value = '  -0.27%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = '1.190'
$ws.Range("E49").Value = '  +1.71%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '1.929'
$ws.Range("E50").Value = '  -2.64%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.06820'
$ws.Range("E51").Value = '  -1.67%  '

# Remove the temporary text-number-format so styling matches the original
# (cells keep their string value/type, just without an explicit style index).
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).ClearFormats()
}
